$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.896.82"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.815.96"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.14"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.51"
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +8.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.33"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.05"
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.83"
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("D15").Value = "3.256.99"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "2.818.17"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "51.870.43"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.68"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("E20").Value = "  +3.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.71"
$ws.Range("E21").Value = "  +4.67%  "
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.54"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.23"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.44"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.10"
$ws.Range("E30").Value = "  +9.96%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.22"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.44"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("E34").Value = "  +10.94%  "
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0872"
$ws.Range("E36").Value = "  +4.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.82"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.51"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.14"
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.30"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  +7.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.40"
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("D48").Value = "2.109.60"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.935"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  +10.23%  "
$ws.Range("E51").Value = "  -1.59%  "
